# Weekly data refresh: two new daily price records for
# "Vega Modelo de Temuco - Apio" are appended to the historical log.
# Because the log isn't stored in date order, the new rows land in the
# middle of the existing block (at row 382 and at row 506), pushing the
# following rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record #1 at row 382 (date serial 45120 = 2023-07-13) ---
$ws.Rows.Item(382).Insert()

$ws.Range("A382").Value = 10
$ws.Range("B382").Value = "Vega Modelo de Temuco"
$ws.Range("C382").Value = "La Araucanía"
$ws.Range("D382").Value = 45120
$ws.Range("E382").Value = 9
$ws.Range("F382").Value = 100112017
$ws.Range("G382").Value = "Apio"
$ws.Range("H382").Value = "Americana (o)"
$ws.Range("I382").Value = "Primera"
$ws.Range("J382").Value = 300
$ws.Range("K382").Value = 7000
$ws.Range("L382").Value = 8000
$ws.Range("M382").Value = 7500
$ws.Range("N382").Value = "`$/docena de matas"
$ws.Range("O382").Value = "Provincia de Limarí"
$ws.Range("P382").Value = 1250
$ws.Range("Q382").Value = 6
$ws.Range("R382").Value = "Hortaliza"

# --- Insert new record #2 at row 506 (date serial 45121 = 2023-07-14) ---
$ws.Rows.Item(506).Insert()

$ws.Range("A506").Value = 10
$ws.Range("B506").Value = "Vega Modelo de Temuco"
$ws.Range("C506").Value = "La Araucanía"
$ws.Range("D506").Value = 45121
$ws.Range("E506").Value = 9
$ws.Range("F506").Value = 100112017
$ws.Range("G506").Value = "Apio"
$ws.Range("H506").Value = "Americana (o)"
$ws.Range("I506").Value = "Primera"
$ws.Range("J506").Value = 40
$ws.Range("K506").Value = 7000
$ws.Range("L506").Value = 7000
$ws.Range("M506").Value = 7000
$ws.Range("N506").Value = "`$/docena de matas"
$ws.Range("O506").Value = "Provincia del Elquí"
$ws.Range("P506").Value = 1167
$ws.Range("Q506").Value = 6
$ws.Range("R506").Value = "Hortaliza"
